# Remove the "popdens" column (column D) from the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:D21").EntireColumn.Delete()

# Move selection, matching the author's resulting state.
$ws.Range("H5").Select()
